$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.528.81"
$ws.Range("E2").Value = "  -0.37%  "

$ws.Range("D3").Value = "3.314.83"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.14%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "576.30"
$ws.Range("E5").Value = "  -0.33%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.89"
$ws.Range("E6").Value = "  -3.57%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.11%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.590"
$ws.Range("E8").Value = "  +0.95%  "

$ws.Range("D9").Value = "3.313.53"
$ws.Range("E9").Value = "  +1.11%  "

$ws.Range("E10").Value = "  +0.97%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.577"
$ws.Range("E11").Value = "  +0.35%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "45.68"
$ws.Range("E12").Value = "  +0.08%  "

$ws.Range("E13").Value = "  -0.30%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "711.72"
$ws.Range("E14").Value = "  +3.00%  "

$ws.Range("D15").Value = "3.860.25"
$ws.Range("E15").Value = "  +1.28%  "

$ws.Range("E16").Value = "  +0.24%  "

$ws.Range("D17").Value = "67.514.89"
$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("E18").Value = "  -1.01%  "

$ws.Range("D19").Value = "3.328.56"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("E20").Value = "  -0.37%  "

$ws.Range("E21").Value = "  +1.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.887"
$ws.Range("E22").Value = "  -0.50%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.36"
$ws.Range("E23").Value = "  +4.00%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.82"
$ws.Range("E24").Value = "  -1.51%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "98.58"
$ws.Range("E25").Value = "  +0.30%  "

$ws.Range("E26").Value = "  -2.45%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.66"
$ws.Range("E27").Value = "  -2.79%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.30"

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.42"
$ws.Range("E29").Value = "  +2.25%  "

$ws.Range("E30").Value = "  +0.58%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.09"
$ws.Range("E31").Value = "  +5.62%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "569.74"
$ws.Range("E32").Value = "  -2.52%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "10.94"
$ws.Range("E33").Value = "  +0.78%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.104"
$ws.Range("E34").Value = "  +0.64%  "

$ws.Range("B35").Value = "Dai"
$ws.Range("C35").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("E35").Value = "  +0.17%  "

$ws.Range("B36").Value = "Maker"
$ws.Range("C36").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D36").Value = "3.698.19"
$ws.Range("E36").Value = "  -4.75%  "

$ws.Range("B37").Value = "OKB"
$ws.Range("C37").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "56.87"
$ws.Range("E37").Value = "  +2.71%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  -2.79%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "34.15"
$ws.Range("E39").Value = "  +5.81%  "

$ws.Range("E40").Value = "  -0.56%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.15"
$ws.Range("E41").Value = "  -2.81%  "

$ws.Range("E42").Value = "  -1.11%  "

$ws.Range("B43").Value = "TheGraph"
$ws.Range("C43").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.332"
$ws.Range("E43").Value = "  +0.33%  "

$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.28"
$ws.Range("E44").Value = "  -3.49%  "

$ws.Range("D45").Value = "0.0₃0663"
$ws.Range("E45").Value = "  -1.75%  "

$ws.Range("E46").Value = "  -2.12%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.67"
$ws.Range("E47").Value = "  +6.88%  "

$ws.Range("E48").Value = "  +0.06%  "

$ws.Range("E49").Value = "  -0.36%  "

$ws.Range("E50").Value = "  -5.35%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "128.63"
$ws.Range("E51").Value = "  -1.06%  "
